$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A30").Value = 45979
$ws.Range("A30").NumberFormat = "d-mmm-yy"
$ws.Range("B30").Value = 5610
$ws.Range("C30").Value = 4382
$ws.Range("D30").Value = 3898
$ws.Range("E30").Value = 315
$ws.Range("F30").Value = 97
$ws.Range("G30").Value = 66
$ws.Range("H30").Value = 6
$ws.Range("I30").Value = 0

$ws.Range("A30:I30").Select() | Out-Null
